# Apply "Scenarios with non-universal participation, AR6 data" edit.
# - Adds two new scenario columns (Dem USA, Non-Dem USA) plus a combined USA
#   column at the end of the country panel, shifting the tail of the country
#   list (Uzbekistan..Zimbabwe) left by one column.
# - Replaces the row-2 statistic label and refreshes all of row 2's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @(
    'Afghanistan',
    'Angola',
    'Albania',
    'Andorra',
    'United Arab Emirates',
    'Argentina',
    'Armenia',
    'Antigua and Barbuda',
    'Australia',
    'Austria',
    'Azerbaijan',
    'Burundi',
    'Belgium',
    'Benin',
    'Burkina Faso',
    'Bangladesh',
    'Bulgaria',
    'Bahrain',
    'Bahamas',
    'Bosnia and Herzegovina',
    'Belarus',
    'Belize',
    'Bolivia',
    'Brazil',
    'Barbados',
    'Brunei',
    'Bhutan',
    'Botswana',
    'Central African Republic',
    'Canada',
    'Switzerland',
    'Chile',
    'China',
    'Ivory Coast',
    'Cameroon',
    'Democratic Republic of the Congo',
    'Republic of Congo',
    'Colombia',
    'Comoros',
    'Cape Verde',
    'Costa Rica',
    'Cuba',
    'Cyprus',
    'Czech Republic',
    'Germany',
    'Djibouti',
    'Dominica',
    'Denmark',
    'Dominican Republic',
    'Algeria',
    'Ecuador',
    'Egypt',
    'Eritrea',
    'Spain',
    'Estonia',
    'Ethiopia',
    'Finland',
    'Fiji',
    'France',
    'Micronesia (country)',
    'Gabon',
    'UK',
    'Georgia',
    'Ghana',
    'Guinea',
    'Gambia',
    'Guinea-Bissau',
    'Equatorial Guinea',
    'Greece',
    'Grenada',
    'Guatemala',
    'Guyana',
    'Honduras',
    'Croatia',
    'Haiti',
    'Hungary',
    'Indonesia',
    'India',
    'Ireland',
    'Iran',
    'Iraq',
    'Iceland',
    'Israel',
    'Italy',
    'Jamaica',
    'Jordan',
    'Japan',
    'Kazakhstan',
    'Kenya',
    'Kyrgyzstan',
    'Cambodia',
    'Kiribati',
    'Saint Kitts and Nevis',
    'South Korea',
    'Kuwait',
    'Laos',
    'Lebanon',
    'Liberia',
    'Libya',
    'Saint Lucia',
    'Liechtenstein',
    'Sri Lanka',
    'Lesotho',
    'Lithuania',
    'Luxembourg',
    'Latvia',
    'Morocco',
    'Moldova',
    'Madagascar',
    'Maldives',
    'Mexico',
    'Marshall Islands',
    'North Macedonia',
    'Mali',
    'Malta',
    'Myanmar',
    'Montenegro',
    'Mongolia',
    'Mozambique',
    'Mauritania',
    'Mauritius',
    'Malawi',
    'Malaysia',
    'Namibia',
    'Niger',
    'Nigeria',
    'Nicaragua',
    'Netherlands',
    'Norway',
    'Nepal',
    'Nauru',
    'New Zealand',
    'Oman',
    'Pakistan',
    'Panama',
    'Peru',
    'Philippines',
    'Palau',
    'Papua New Guinea',
    'Poland',
    'North Korea',
    'Portugal',
    'Paraguay',
    'Palestine',
    'Qatar',
    'Romania',
    'Russia',
    'Rwanda',
    'Saudi Arabia',
    'Sudan',
    'Senegal',
    'Singapore',
    'Solomon Islands',
    'Sierra Leone',
    'El Salvador',
    'Somalia',
    'Serbia',
    'South Sudan',
    'Sao Tome and Principe',
    'Suriname',
    'Slovakia',
    'Slovenia',
    'Sweden',
    'Eswatini',
    'Seychelles',
    'Syria',
    'Chad',
    'Togo',
    'Thailand',
    'Tajikistan',
    'Turkmenistan',
    'Timor',
    'Tonga',
    'Trinidad and Tobago',
    'Tunisia',
    'Turkey',
    'Tuvalu',
    'Tanzania',
    'Uganda',
    'Ukraine',
    'Uruguay',
    'Uzbekistan',
    'Saint Vincent and the Grenadines',
    'Venezuela',
    'Vietnam',
    'Vanuatu',
    'Samoa',
    'Yemen',
    'South Africa',
    'Zambia',
    'Zimbabwe',
    'Dem USA',
    'Non-Dem USA',
    'USA'
)

$values = @(
    0.442607759855407,
    0.0363536571999655,
    0.034703609775288,
    $null,
    $null,
    0,
    0.0256863986529112,
    $null,
    $null,
    -0.0116626134644162,
    0,
    1.4709087207359,
    -0.0306372508031646,
    0.12104285292686,
    0.397682602404877,
    0.122900124175341,
    0,
    $null,
    $null,
    0,
    0,
    0,
    0,
    0,
    $null,
    $null,
    0.0142264644065338,
    0,
    1.06289942058061,
    -0.0218098754504761,
    -0.0109297623983196,
    $null,
    0,
    0.0699447047896936,
    0.166924533875106,
    0.83598408692024,
    0,
    0,
    0.184039347608431,
    0.00967565394739187,
    0,
    0.000378398450013442,
    -0.00642636942332954,
    -0.0339079420262336,
    -0.0125275397455969,
    0.10046898710192,
    0,
    -0.00629443676154853,
    0,
    0,
    0,
    0.00611945496414839,
    0.441877554415933,
    -0.00520198258476762,
    -0.0377241945038636,
    0.435619187319076,
    -0.0112525711859789,
    0.0285925010970675,
    -0.00549881143792919,
    0.0600040226405731,
    0,
    -0.00722320759664687,
    0.00920144043933313,
    0.0613071884669687,
    0.206697427830322,
    0.466268485014821,
    0.556466144145017,
    0,
    -0.00240428666241817,
    0,
    0.011063572316679,
    0,
    0.0732534504730981,
    -0.00770412686575331,
    0.239818813414395,
    -0.0159233152086035,
    0.0107923373320534,
    0.0269657637302264,
    -0.00843252236211339,
    0,
    0,
    $null,
    $null,
    -0.00801015730541485,
    0,
    0,
    -0.00703889615832972,
    $null,
    0.110774487067026,
    0.0221046114216444,
    0.113640404369409,
    0.199528461102441,
    $null,
    -0.0359451846669731,
    $null,
    0.0171256720666182,
    0,
    0.50440643242093,
    0,
    0,
    $null,
    0.0162633852633878,
    0,
    -0.0237662133450578,
    -0.00989873175103032,
    -0.0280415130012458,
    0.04349375032316,
    0.0540543683312584,
    0.702493908685207,
    $null,
    0,
    0,
    0,
    0.424785877035098,
    -0.0514504791319278,
    0.177521122890708,
    0,
    0,
    0.345603258832404,
    0.012554352502925,
    $null,
    0.912860825052543,
    $null,
    0,
    0.714849617938332,
    0.0582353547092945,
    0.082041159510952,
    -0.0106723678437547,
    -0.00645138847787756,
    0.170288272515342,
    0,
    -0.00403775267183949,
    $null,
    0.0965625495302393,
    $null,
    0,
    0.0385582353552717,
    $null,
    0.093229198535131,
    -0.0233160180060687,
    0.427363521268592,
    -0.00237631977058196,
    0,
    0.103922639150555,
    $null,
    0,
    $null,
    0.397774447916352,
    $null,
    0.196251455354028,
    0,
    $null,
    0.133783254435204,
    0.572657174465654,
    0.0151096473898082,
    0.951619693117974,
    0,
    0.97992583675477,
    0.0867064489873992,
    0,
    -0.0210040337213707,
    -0.0254741140578765,
    -0.00448931863020271,
    0.00613231900291519,
    $null,
    0.274290325231816,
    0.561155004743766,
    0,
    0,
    0.149096837497506,
    $null,
    0.184878674715778,
    0.0291827994614112,
    $null,
    0.0238281065389928,
    0,
    0.0645946780307104,
    0.269052124065286,
    0.360758420250027,
    0,
    $null,
    0,
    0,
    0.0972200902726076,
    0.0733079630670761,
    0.0984667030473454,
    0.0374017235368603,
    0.58046366968681,
    0,
    0.145663151967098,
    0.0473984698509576,
    -0.00990695992677429,
    $null,
    $null
)

$colCount = 194

$headerArr = New-Object 'object[,]' 1,$colCount
$valueArr = New-Object 'object[,]' 1,$colCount
for ($i = 0; $i -lt $colCount; $i++) {
    $headerArr[0,$i] = $headers[$i]
    $valueArr[0,$i] = $values[$i]
}

$ws.Range("B1:GM1").Value = $headerArr
$ws.Range("A2").Value = 'Soptimistic_gain_adj_over_gdp_2030'
$ws.Range("B2:GM2").Value = $valueArr

Write-Output "applied edit: headers + row2 updated"
